$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.738.03"
$ws.Range("E2").Value = "'  +1.02%  "

$ws.Range("D3").Value = "'1.659.52"
$ws.Range("E3").Value = "'  +1.12%  "

$ws.Range("E4").Value = "'  +0.04%  "

$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "'  +0.06%  "

$ws.Range("D6").Value = "'302.98"
$ws.Range("E6").Value = "'  -0.26%  "

$ws.Range("D7").Value = "'0.3817"
$ws.Range("E7").Value = "'  +0.75%  "

$ws.Range("D8").Value = "'0.3614"
$ws.Range("E8").Value = "'  -0.32%  "

$ws.Range("D9").Value = "'51.27"
$ws.Range("E9").Value = "'  -1.25%  "

$ws.Range("E10").Value = "'  +0.09%  "

$ws.Range("D11").Value = "'1.230"
$ws.Range("E11").Value = "'  -0.51%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "'  +0.12%  "

$ws.Range("D13").Value = "'22.55"
$ws.Range("E13").Value = "'  -0.17%  "

$ws.Range("D14").Value = "'6.473"
$ws.Range("E14").Value = "'  +0.01%  "

$ws.Range("D15").Value = "'7.414"
$ws.Range("E15").Value = "'  +0.48%  "

$ws.Range("E16").Value = "'  -1.08%  "

$ws.Range("D17").Value = "'1.654.70"
$ws.Range("E17").Value = "'  +1.40%  "

$ws.Range("D18").Value = "'97.91"
$ws.Range("E18").Value = "'  +2.85%  "

$ws.Range("E19").Value = "'  +1.08%  "

$ws.Range("E20").Value = "'  +3.75%  "

$ws.Range("D21").Value = "'17.64"
$ws.Range("E21").Value = "'  +0.45%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "'  +0.05%  "

$ws.Range("D23").Value = "'12.85"
$ws.Range("E23").Value = "'  +2.54%  "

$ws.Range("D24").Value = "'23.726.63"
$ws.Range("E24").Value = "'  +1.01%  "

$ws.Range("D25").Value = "'2.505"
$ws.Range("E25").Value = "'  -0.21%  "

$ws.Range("D26").Value = "'2.995"
$ws.Range("E26").Value = "'  -3.15%  "

$ws.Range("D27").Value = "'21.22"
$ws.Range("E27").Value = "'  -0.01%  "

$ws.Range("D28").Value = "'152.83"
$ws.Range("E28").Value = "'  +0.24%  "

$ws.Range("D29").Value = "'5.225"
$ws.Range("E29").Value = "'  -0.58%  "

$ws.Range("D30").Value = "'134.41"
$ws.Range("E30").Value = "'  +0.73%  "

$ws.Range("D31").Value = "'1.836.22"
$ws.Range("E31").Value = "'  +1.06%  "

$ws.Range("D32").Value = "'7.111"
$ws.Range("E32").Value = "'  +7.44%  "

$ws.Range("D33").Value = "'2.227"
$ws.Range("E33").Value = "'  +3.73%  "

$ws.Range("D34").Value = "'12.14"
$ws.Range("E34").Value = "'  +5.70%  "

$ws.Range("D35").Value = "'1.055"
$ws.Range("E35").Value = "'  -4.11%  "

$ws.Range("D36").Value = "'0.02813"
$ws.Range("E36").Value = "'  +1.73%  "

$ws.Range("D37").Value = "'0.2515"
$ws.Range("E37").Value = "'  +0.43%  "

$ws.Range("D38").Value = "'0.08805"
$ws.Range("E38").Value = "'  +0.31%  "

$ws.Range("D39").Value = "'6.093"
$ws.Range("E39").Value = "'  +0.94%  "

$ws.Range("D40").Value = "'0.07015"
$ws.Range("E40").Value = "'  -0.99%  "

$ws.Range("D41").Value = "'13.01"
$ws.Range("E41").Value = "'  +5.19%  "

$ws.Range("D42").Value = "'0.7000"
$ws.Range("E42").Value = "'  -0.85%  "

$ws.Range("D43").Value = "'1.333"
$ws.Range("E43").Value = "'  -1.62%  "

$ws.Range("D44").Value = "'16.02"
$ws.Range("E44").Value = "'  +1.88%  "

$ws.Range("D45").Value = "'0.6517"
$ws.Range("E45").Value = "'  -0.47%  "

$ws.Range("B46").Value = "'Frax"
$ws.Range("C46").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'1.0000"
$ws.Range("E46").Value = "'  +0.04%  "

$ws.Range("B47").Value = "'NEARProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.309"
$ws.Range("E47").Value = "'  +1.00%  "

$ws.Range("D48").Value = "'3.965"
$ws.Range("E48").Value = "'  -0.10%  "

$ws.Range("D49").Value = "'0.07917"
$ws.Range("E49").Value = "'  -0.88%  "

$ws.Range("D50").Value = "'128.23"
$ws.Range("E50").Value = "'  -0.62%  "

$ws.Range("D51").Value = "'1.183"
$ws.Range("E51").Value = "'  -1.14%  "
